# Shuffle rows 2-27: keep columns A-E as-is, set column F (trial_total) equal
# to column E (trial_block), and re-assign columns G-S (target_cat .. p_perceptual)
# from a permuted source row, per the committed data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 27

# For each destination row (index = row - firstRow), the 1-based worksheet
# row whose G:S content should be copied into it.
$sourceRowForDest = @{
    2  = 6
    3  = 12
    4  = 7
    5  = 8
    6  = 26
    7  = 14
    8  = 23
    9  = 11
    10 = 21
    11 = 19
    12 = 9
    13 = 17
    14 = 4
    15 = 2
    16 = 18
    17 = 13
    18 = 3
    19 = 16
    20 = 22
    21 = 15
    22 = 25
    23 = 20
    24 = 27
    25 = 5
    26 = 24
    27 = 10
}

# Snapshot the original G:S values (columns 7-19) for every row before any
# writes happen, so overlapping reads/writes during the shuffle don't clobber
# data we still need to read later.
$original = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @{}
    for ($c = 7; $c -le 19; $c++) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value()
    }
    $original[$r] = $rowVals
}

for ($r = $firstRow; $r -le $lastRow; $r++) {
    # Column F (trial_total) becomes the same sequential number as column E
    # (trial_block).
    $ws.Cells.Item($r, 6).Value = $ws.Cells.Item($r, 5).Value()

    $src = $sourceRowForDest[$r]
    $srcVals = $original[$src]
    for ($c = 7; $c -le 19; $c++) {
        $ws.Cells.Item($r, $c).Value = $srcVals[$c]
    }
}
